$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 95.083336  # H8: 31 -> 95.083336
$ws.Cells.Item(8, 9).Value = 94.63636  # I8: 24.1 -> 94.63636
$ws.Cells.Item(8, 11).Value = 283.90908  # K8: 72.30000000000001 -> 283.90908
$ws.Cells.Item(8, 13).Value = -144.90908  # M8: 66.69999999999999 -> -144.90908

$ws.Cells.Item(29, 8).Value = 1900  # H29: 2300 -> 1900
$ws.Cells.Item(29, 9).Value = 375  # I29: 400 -> 375
$ws.Cells.Item(29, 11).Value = 1125  # K29: 1200 -> 1125
$ws.Cells.Item(29, 13).Value = -844  # M29: -919 -> -844

$ws.Cells.Item(62, 8).Value = 6750  # H62: 5833.3335 -> 6750
$ws.Cells.Item(62, 9).Value = 3500  # I62: 3750 -> 3500
$ws.Cells.Item(62, 11).Value = 3500  # K62: 3750 -> 3500
$ws.Cells.Item(62, 13).Value = -2876  # M62: -3126 -> -2876

$ws.Cells.Item(65, 8).Value = 6750  # H65: 5833.3335 -> 6750
$ws.Cells.Item(65, 9).Value = 3500  # I65: 3750 -> 3500
$ws.Cells.Item(65, 11).Value = 17500  # K65: 18750 -> 17500
$ws.Cells.Item(65, 13).Value = -14380  # M65: -15630 -> -14380

$ws.Cells.Item(141, 8).Value = 20414.5  # H141: 23697.6 -> 20414.5
$ws.Cells.Item(141, 9).Value = 20414.5  # I141: 23697.6 -> 20414.5
$ws.Cells.Item(141, 11).Value = 61243.5  # K141: 71092.79999999999 -> 61243.5
$ws.Cells.Item(141, 13).Value = -56063.5  # M141: -65912.79999999999 -> -56063.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(6, 8).Value = 5000  # H6: 0 -> 5000
$ws.Cells.Item(6, 10).Value = 5000  # J6: 0 -> 5000
$ws.Cells.Item(6, 12).Value = 5000  # L6: 0 -> 5000
$ws.Cells.Item(6, 14).Value = -5346  # N6: <<MISSING>> -> -5346

$ws.Cells.Item(15, 8).Value = 0  # H15: 20000 -> 0
$ws.Cells.Item(15, 10).Value = 0  # J15: 20000 -> 0
$ws.Cells.Item(15, 12).Value = 0  # L15: 20000 -> 0
$ws.Cells.Item(15, 14).ClearContents()  # N15: -20700 -> (deleted)

$ws.Cells.Item(54, 8).Value = 19000  # H54: 50000 -> 19000
$ws.Cells.Item(54, 10).Value = 19000  # J54: 50000 -> 19000
$ws.Cells.Item(54, 12).Value = 19000  # L54: 50000 -> 19000
$ws.Cells.Item(54, 14).Value = -20538  # N54: -51538 -> -20538

$ws.Cells.Item(92, 8).Value = 49999.5  # H92: 62999.4 -> 49999.5
$ws.Cells.Item(92, 9).Value = 0  # I92: 90000 -> 0
$ws.Cells.Item(92, 10).Value = 49999.5  # J92: 56249.25 -> 49999.5
$ws.Cells.Item(92, 11).Value = 0  # K92: 90000 -> 0
$ws.Cells.Item(92, 12).Value = 49999.5  # L92: 56249.25 -> 49999.5
$ws.Cells.Item(92, 13).ClearContents()  # M92: -87504 -> (deleted)
$ws.Cells.Item(92, 14).Value = -54991.5  # N92: -61241.25 -> -54991.5

$ws.Cells.Item(95, 8).Value = 20470.5  # H95: 20964.6 -> 20470.5
$ws.Cells.Item(95, 10).Value = 20470.5  # J95: 20964.6 -> 20470.5
$ws.Cells.Item(95, 12).Value = 20470.5  # L95: 20964.6 -> 20470.5
$ws.Cells.Item(95, 14).Value = -25962.5  # N95: -26456.6 -> -25962.5

$ws.Cells.Item(96, 8).Value = 35114  # H96: 35172 -> 35114
$ws.Cells.Item(96, 10).Value = 35114  # J96: 35172 -> 35114
$ws.Cells.Item(96, 12).Value = 35114  # L96: 35172 -> 35114
$ws.Cells.Item(96, 14).Value = -40606  # N96: -40664 -> -40606

$ws.Cells.Item(106, 8).Value = 32546  # H106: 31645.5 -> 32546
$ws.Cells.Item(106, 10).Value = 32546  # J106: 31645.5 -> 32546
$ws.Cells.Item(106, 12).Value = 32546  # L106: 31645.5 -> 32546
$ws.Cells.Item(106, 14).Value = -35070  # N106: -34169.5 -> -35070

$ws.Cells.Item(112, 8).Value = 30000  # H112: 0 -> 30000
$ws.Cells.Item(112, 10).Value = 30000  # J112: 0 -> 30000
$ws.Cells.Item(112, 12).Value = 30000  # L112: 0 -> 30000
$ws.Cells.Item(112, 14).Value = -32954  # N112: <<MISSING>> -> -32954

$ws.Cells.Item(122, 8).Value = 0  # H122: 1700 -> 0
$ws.Cells.Item(122, 9).Value = 0  # I122: 1700 -> 0
$ws.Cells.Item(122, 11).Value = 0  # K122: 5100 -> 0
$ws.Cells.Item(122, 13).ClearContents()  # M122: -2650 -> (deleted)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(11, 8).Value = 866.6667  # H11: 570 -> 866.6667
$ws.Cells.Item(11, 9).Value = 300  # I11: 212.5 -> 300
$ws.Cells.Item(11, 11).Value = 300  # K11: 212.5 -> 300
$ws.Cells.Item(11, 13).Value = -160  # M11: -72.5 -> -160

$ws.Cells.Item(20, 8).Value = 1953.1428  # H20: 2160.6667 -> 1953.1428
$ws.Cells.Item(20, 9).Value = 1536.4  # I20: 1743.5 -> 1536.4
$ws.Cells.Item(20, 11).Value = 1536.4  # K20: 1743.5 -> 1536.4
$ws.Cells.Item(20, 13).Value = -1289.4  # M20: -1496.5 -> -1289.4

$ws.Cells.Item(86, 8).Value = 1321.4  # H86: 1434.8889 -> 1321.4
$ws.Cells.Item(86, 9).Value = 1351.875  # I86: 1502.1428 -> 1351.875
$ws.Cells.Item(86, 11).Value = 1351.875  # K86: 1502.1428 -> 1351.875
$ws.Cells.Item(86, 13).Value = -228.875  # M86: -379.1428000000001 -> -228.875

$ws.Cells.Item(89, 8).Value = 1321.4  # H89: 1434.8889 -> 1321.4
$ws.Cells.Item(89, 9).Value = 1351.875  # I89: 1502.1428 -> 1351.875
$ws.Cells.Item(89, 11).Value = 6759.375  # K89: 7510.714 -> 6759.375
$ws.Cells.Item(89, 13).Value = -1143.375  # M89: -1894.714 -> -1143.375

$ws.Cells.Item(94, 8).Value = 1374.3334  # H94: 1483.75 -> 1374.3334
$ws.Cells.Item(94, 9).Value = 1100  # I94: 1250.25 -> 1100
$ws.Cells.Item(94, 11).Value = 1100  # K94: 1250.25 -> 1100
$ws.Cells.Item(94, 13).Value = -649  # M94: -799.25 -> -649

$ws.Cells.Item(132, 8).Value = 0  # H132: 48000 -> 0
$ws.Cells.Item(132, 10).Value = 0  # J132: 48000 -> 0
$ws.Cells.Item(132, 12).Value = 0  # L132: 48000 -> 0
$ws.Cells.Item(132, 14).ClearContents()  # N132: -58120 -> (deleted)

$ws.Cells.Item(140, 8).Value = 95779.11  # H140: 95779.125 -> 95779.11
$ws.Cells.Item(140, 10).Value = 95779.11  # J140: 95779.125 -> 95779.11
$ws.Cells.Item(140, 12).Value = 95779.11  # L140: 95779.125 -> 95779.11
$ws.Cells.Item(140, 14).Value = -106139.11  # N140: -106139.125 -> -106139.11

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 7335.6665  # H86: 7003.5 -> 7335.6665
$ws.Cells.Item(86, 10).Value = 8000  # J86: 0 -> 8000
$ws.Cells.Item(86, 12).Value = 8000  # L86: 0 -> 8000
$ws.Cells.Item(86, 14).Value = -10246  # N86: <<MISSING>> -> -10246

$ws.Cells.Item(89, 8).Value = 7335.6665  # H89: 7003.5 -> 7335.6665
$ws.Cells.Item(89, 10).Value = 8000  # J89: 0 -> 8000
$ws.Cells.Item(89, 12).Value = 40000  # L89: 0 -> 40000
$ws.Cells.Item(89, 14).Value = -51232  # N89: <<MISSING>> -> -51232

$ws.Cells.Item(122, 8).Value = 1325  # H122: 1625 -> 1325
$ws.Cells.Item(122, 9).Value = 1325  # I122: 1625 -> 1325
$ws.Cells.Item(122, 11).Value = 3975  # K122: 4875 -> 3975
$ws.Cells.Item(122, 13).Value = -1525  # M122: -2425 -> -1525

$ws.Cells.Item(141, 8).Value = 1341666  # H141: 1999999 -> 1341666
$ws.Cells.Item(141, 10).Value = 1341666  # J141: 1999999 -> 1341666
$ws.Cells.Item(141, 12).Value = 1341666  # L141: 1999999 -> 1341666
$ws.Cells.Item(141, 14).Value = -1352026  # N141: -2010359 -> -1352026

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 6066.6665  # H55: 200 -> 6066.6665
$ws.Cells.Item(55, 10).Value = 9000  # J55: 0 -> 9000
$ws.Cells.Item(55, 12).Value = 27000  # L55: 0 -> 27000
$ws.Cells.Item(55, 14).Value = -27354  # N55: <<MISSING>> -> -27354

$ws.Cells.Item(133, 8).Value = 0  # H133: 3025 -> 0
$ws.Cells.Item(133, 9).Value = 0  # I133: 3025 -> 0
$ws.Cells.Item(133, 11).Value = 0  # K133: 9075 -> 0
$ws.Cells.Item(133, 13).ClearContents()  # M133: -4015 -> (deleted)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 2004950  # H7: 10000 -> 2004950
$ws.Cells.Item(7, 9).Value = 0  # I7: 10000 -> 0
$ws.Cells.Item(7, 10).Value = 2004950  # J7: 0 -> 2004950
$ws.Cells.Item(7, 11).Value = 0  # K7: 10000 -> 0
$ws.Cells.Item(7, 12).Value = 2004950  # L7: 0 -> 2004950
$ws.Cells.Item(7, 13).ClearContents()  # M7: -9888 -> (deleted)
$ws.Cells.Item(7, 14).Value = -2005174  # N7: <<MISSING>> -> -2005174

$ws.Cells.Item(8, 8).Value = 2004950  # H8: 10000 -> 2004950
$ws.Cells.Item(8, 9).Value = 0  # I8: 10000 -> 0
$ws.Cells.Item(8, 10).Value = 2004950  # J8: 0 -> 2004950
$ws.Cells.Item(8, 11).Value = 0  # K8: 10000 -> 0
$ws.Cells.Item(8, 12).Value = 2004950  # L8: 0 -> 2004950
$ws.Cells.Item(8, 13).ClearContents()  # M8: -9861 -> (deleted)
$ws.Cells.Item(8, 14).Value = -2005228  # N8: <<MISSING>> -> -2005228

$ws.Cells.Item(22, 8).Value = 5666.6665  # H22: 1500 -> 5666.6665
$ws.Cells.Item(22, 9).Value = 1000  # I22: 0 -> 1000
$ws.Cells.Item(22, 10).Value = 15000  # J22: 1500 -> 15000
$ws.Cells.Item(22, 11).Value = 1000  # K22: 0 -> 1000
$ws.Cells.Item(22, 12).Value = 15000  # L22: 1500 -> 15000
$ws.Cells.Item(22, 13).Value = -471  # M22: <<MISSING>> -> -471
$ws.Cells.Item(22, 14).Value = -16058  # N22: -2558 -> -16058

$ws.Cells.Item(25, 8).Value = 2000  # H25: 7766.6665 -> 2000
$ws.Cells.Item(25, 10).Value = 2000  # J25: 7766.6665 -> 2000
$ws.Cells.Item(25, 12).Value = 2000  # L25: 7766.6665 -> 2000
$ws.Cells.Item(25, 14).Value = -3058  # N25: -8824.666499999999 -> -3058

$ws.Cells.Item(70, 8).Value = 4997.5  # H70: 5333.3335 -> 4997.5
$ws.Cells.Item(70, 9).Value = 4496.6665  # I70: 4750 -> 4496.6665
$ws.Cells.Item(70, 11).Value = 4496.6665  # K70: 4750 -> 4496.6665
$ws.Cells.Item(70, 13).Value = -4226.6665  # M70: -4480 -> -4226.6665

$ws.Cells.Item(73, 8).Value = 4997.5  # H73: 5333.3335 -> 4997.5
$ws.Cells.Item(73, 9).Value = 4496.6665  # I73: 4750 -> 4496.6665
$ws.Cells.Item(73, 11).Value = 4496.6665  # K73: 4750 -> 4496.6665
$ws.Cells.Item(73, 13).Value = -3560.6665  # M73: -3814 -> -3560.6665

$ws.Cells.Item(102, 8).Value = 493  # H102: 332.1111 -> 493
$ws.Cells.Item(102, 9).Value = 493  # I102: 336.125 -> 493
$ws.Cells.Item(102, 10).Value = 0  # J102: 300 -> 0
$ws.Cells.Item(102, 11).Value = 493  # K102: 336.125 -> 493
$ws.Cells.Item(102, 12).Value = 0  # L102: 300 -> 0
$ws.Cells.Item(102, 13).Value = 1129  # M102: 1285.875 -> 1129
$ws.Cells.Item(102, 14).ClearContents()  # N102: -3544 -> (deleted)

$ws.Cells.Item(134, 8).Value = 40000  # H134: 90000 -> 40000
$ws.Cells.Item(134, 10).Value = 40000  # J134: 90000 -> 40000
$ws.Cells.Item(134, 12).Value = 120000  # L134: 270000 -> 120000
$ws.Cells.Item(134, 14).Value = -125070  # N134: -275070 -> -125070

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(50, 8).Value = 0  # H50: 54084 -> 0
$ws.Cells.Item(50, 10).Value = 0  # J50: 54084 -> 0
$ws.Cells.Item(50, 12).Value = 0  # L50: 54084 -> 0
$ws.Cells.Item(50, 14).ClearContents()  # N50: -55358 -> (deleted)

$ws.Cells.Item(54, 8).Value = 49084  # H54: 39542 -> 49084
$ws.Cells.Item(54, 10).Value = 49084  # J54: 39542 -> 49084
$ws.Cells.Item(54, 12).Value = 49084  # L54: 39542 -> 49084
$ws.Cells.Item(54, 14).Value = -50372  # N54: -40830 -> -50372

$ws.Cells.Item(82, 8).Value = 500  # H82: 525 -> 500
$ws.Cells.Item(82, 9).Value = 0  # I82: 500 -> 0
$ws.Cells.Item(82, 10).Value = 500  # J82: 550 -> 500
$ws.Cells.Item(82, 11).Value = 0  # K82: 500 -> 0
$ws.Cells.Item(82, 12).Value = 500  # L82: 550 -> 500
$ws.Cells.Item(82, 13).ClearContents()  # M82: -139 -> (deleted)
$ws.Cells.Item(82, 14).Value = -1222  # N82: -1272 -> -1222

$ws.Cells.Item(85, 8).Value = 500  # H85: 525 -> 500
$ws.Cells.Item(85, 9).Value = 0  # I85: 500 -> 0
$ws.Cells.Item(85, 10).Value = 500  # J85: 550 -> 500
$ws.Cells.Item(85, 11).Value = 0  # K85: 500 -> 0
$ws.Cells.Item(85, 12).Value = 500  # L85: 550 -> 500
$ws.Cells.Item(85, 13).ClearContents()  # M85: 748 -> (deleted)
$ws.Cells.Item(85, 14).Value = -2996  # N85: -3046 -> -2996

$ws.Cells.Item(93, 8).Value = 1482.4  # H93: 1586.5555 -> 1482.4
$ws.Cells.Item(93, 10).Value = 977.6667  # J93: 1194 -> 977.6667
$ws.Cells.Item(93, 12).Value = 977.6667  # L93: 1194 -> 977.6667
$ws.Cells.Item(93, 14).Value = -3473.6667  # N93: -3690 -> -3473.6667

$ws.Cells.Item(95, 8).Value = 45000  # H95: 0 -> 45000
$ws.Cells.Item(95, 10).Value = 45000  # J95: 0 -> 45000
$ws.Cells.Item(95, 12).Value = 45000  # L95: 0 -> 45000
$ws.Cells.Item(95, 14).Value = -50492  # N95: <<MISSING>> -> -50492

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(13, 8).Value = 2298  # H13: 1599.3334 -> 2298
$ws.Cells.Item(13, 9).Value = 0  # I13: 500 -> 0
$ws.Cells.Item(13, 10).Value = 2298  # J13: 2149 -> 2298
$ws.Cells.Item(13, 11).Value = 0  # K13: 500 -> 0
$ws.Cells.Item(13, 12).Value = 2298  # L13: 2149 -> 2298
$ws.Cells.Item(13, 13).ClearContents()  # M13: -360 -> (deleted)
$ws.Cells.Item(13, 14).Value = -2578  # N13: -2429 -> -2578

$ws.Cells.Item(54, 8).Value = 21666.334  # H54: 9750 -> 21666.334
$ws.Cells.Item(54, 9).Value = 17499.5  # I54: 9750 -> 17499.5
$ws.Cells.Item(54, 10).Value = 30000  # J54: 0 -> 30000
$ws.Cells.Item(54, 11).Value = 17499.5  # K54: 9750 -> 17499.5
$ws.Cells.Item(54, 12).Value = 30000  # L54: 0 -> 30000
$ws.Cells.Item(54, 13).Value = -16979.5  # M54: -9230 -> -16979.5
$ws.Cells.Item(54, 14).Value = -31040  # N54: <<MISSING>> -> -31040

$ws.Cells.Item(122, 8).Value = 25999.75  # H122: 13743.5 -> 25999.75
$ws.Cells.Item(122, 9).Value = 0  # I122: 1487.5 -> 0
$ws.Cells.Item(122, 10).Value = 25999.75  # J122: 25999.5 -> 25999.75
$ws.Cells.Item(122, 11).Value = 0  # K122: 4462.5 -> 0
$ws.Cells.Item(122, 12).Value = 77999.25  # L122: 77998.5 -> 77999.25
$ws.Cells.Item(122, 13).ClearContents()  # M122: -2012.5 -> (deleted)
$ws.Cells.Item(122, 14).Value = -82899.25  # N122: -82898.5 -> -82899.25
